# Add a "Skill Description" column (E) to the PDP sheet and populate a
# full-name "SFIA Level" column (B) per SkillCode, shifting the old
# Keycode/Description values right into C/D/E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each SkillCode (column A) to its full skill name.
$skillNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "SORC"       = "Sourcing"
    "SUPP"       = "Supplier management"
    "MADE"       = "MADE"
}

# Find last used row on the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# New header for column E.
$ws.Cells.Item(1, 5).Value2 = "Skill Description"

for ($r = 2; $r -le $lastRow; $r++) {
    $skillCode = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($skillCode)) { continue }

    $oldKeycode = $ws.Cells.Item($r, 3).Value2      # old Keycode number (C)
    $oldDescription = $ws.Cells.Item($r, 4).Value2  # old Description text (D)

    $fullName = $skillNames[$skillCode]
    if ([string]::IsNullOrEmpty($fullName)) { $fullName = $skillCode }

    $ws.Cells.Item($r, 2).Value2 = $fullName        # B: SFIA Level -> full skill name
    $ws.Cells.Item($r, 3).Value2 = 1                # C: Keycode -> constant 1
    $ws.Cells.Item($r, 4).Value2 = $oldKeycode      # D: Description -> old Keycode number
    $ws.Cells.Item($r, 5).Value2 = $oldDescription  # E: Skill Description -> old Description text
}
